# chore(runtime): publish files + archive (2025-11-04 11:02:35)
# Refresh KHL referee stats tables ("Главные" / "Линейные") with the
# latest scrape: a handful of games/PIM counters moved for a few
# referees, and the as_of_utc timestamp on every data row advances to
# the new scrape time.

$wb = $excel.ActiveWorkbook
$newTimestamp = "2025-11-04 03:02:32"

# ---- "Главные" (head referees) sheet -------------------------------
$ws = $wb.Worksheets.Item("Главные")

# Updated per-referee counters (row -> Games_KHL..PIM_2min_away, plus
# PIM_K_away for a couple of rows)
$ws.Range("C5").Value  = 21
$ws.Range("D5").Value  = 339
$ws.Range("E5").Value  = 183
$ws.Range("F5").Value  = 156
$ws.Range("G5").Value  = 16.14
$ws.Range("H5").Value  = 8.710000000000001
$ws.Range("I5").Value  = 7.43
$ws.Range("J5").Value  = 89
$ws.Range("K5").Value  = 78
$ws.Range("W5").Value  = 16

$ws.Range("C18").Value = 21
$ws.Range("D18").Value = 300
$ws.Range("E18").Value = 140
$ws.Range("F18").Value = 160
$ws.Range("G18").Value = 14.29
$ws.Range("H18").Value = 6.67
$ws.Range("I18").Value = 7.62
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 75

$ws.Range("C20").Value = 20
$ws.Range("D20").Value = 333
$ws.Range("E20").Value = 134
$ws.Range("F20").Value = 199
$ws.Range("G20").Value = 16.65
$ws.Range("H20").Value = 6.7
$ws.Range("I20").Value = 9.949999999999999
$ws.Range("J20").Value = 62
$ws.Range("K20").Value = 72
$ws.Range("W20").Value = 8

$ws.Range("C22").Value = 14
$ws.Range("D22").Value = 284
$ws.Range("E22").Value = 108
$ws.Range("F22").Value = 176
$ws.Range("G22").Value = 20.29
$ws.Range("H22").Value = 7.71
$ws.Range("I22").Value = 12.57
$ws.Range("J22").Value = 54
$ws.Range("K22").Value = 58

# Every data row's scrape timestamp advances
for ($r = 2; $r -le 26; $r++) {
    $ws.Range("AA$r").Value = $newTimestamp
}

# ---- "Линейные" (line referees) sheet -------------------------------
$ws = $wb.Worksheets.Item("Линейные")

$ws.Range("C4").Value  = 10
$ws.Range("D4").Value  = 172
$ws.Range("E4").Value  = 70
$ws.Range("F4").Value  = 102
$ws.Range("G4").Value  = 17.2
$ws.Range("H4").Value  = 7
$ws.Range("I4").Value  = 10.2
$ws.Range("J4").Value  = 35
$ws.Range("K4").Value  = 41
$ws.Range("W4").Value  = 6

$ws.Range("C9").Value  = 20
$ws.Range("D9").Value  = 375
$ws.Range("E9").Value  = 173
$ws.Range("F9").Value  = 202
$ws.Range("G9").Value  = 18.75
$ws.Range("H9").Value  = 8.65
$ws.Range("I9").Value  = 10.1
$ws.Range("J9").Value  = 74
$ws.Range("K9").Value  = 91

$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 486
$ws.Range("E21").Value = 205
$ws.Range("F21").Value = 281
$ws.Range("G21").Value = 21.13
$ws.Range("H21").Value = 8.91
$ws.Range("I21").Value = 12.22
$ws.Range("J21").Value = 95
$ws.Range("K21").Value = 113

# Every data row's scrape timestamp advances
for ($r = 2; $r -le 26; $r++) {
    $ws.Range("AA$r").Value = $newTimestamp
}
